# Applies the "newest changes to bigram model" update:
#  - appends 6 new trigram-pattern rows (9-14) to the sheet, reusing the
#    same column layout/values pattern as existing rows
#  - updates the sheet view (top-left cell + selection) to match the
#    post-edit state

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newLabels = @(
    "DeveloperControlsProcess",
    "DevelopingInPairs",
    "DistributeWorkEvenly",
    "DivideAndConquer",
    "DomainExpertiseInRoles",
    "FeatureAssignment"
)

$values = @(1.01, 1.01, 1.01, 1.01, 1.01, 0.89, 0.89, 0.89, 1.01, 1.01, 1.01)

$startRow = 9
for ($i = 0; $i -lt $newLabels.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newLabels[$i]

    for ($col = 2; $col -le 12; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 2]
    }

    $ws.Cells.Item($row, 13).Value = $row - 1
}

# Match the saved view state: scrolled to column I, selection on M9:M14.
$ws.Range("M9:M14").Select()
$ws.Application.ActiveWindow.ScrollColumn = 9
